# Applies the "Automatic update of files" change to the
# "Avverkningsanmälningar" sheet of Översikt OKÄNT.xlsx:
#   - the "Förändrad" date (column C) for every existing data row (2-32)
#     moves from 2024-11-25 (serial 45621) to 2024-11-27 (serial 45623)
#   - a new case "A 52224-2024" is inserted as row 33, pushing the former
#     row 33 ("A 52574-2024") down to row 34 (its "Förändrad" date is
#     likewise refreshed to 45623)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "Förändrad" column for all existing rows (2 through 32).
for ($r = 2; $r -le 32; $r++) {
  $ws.Cells.Item($r, 3).Value = 45623
}

# Insert a new row above the current last row (33), shifting the old
# row 33 down to row 34 and carrying its formatting (date number format,
# wrap-text style, etc.) onto the newly inserted row 33.
$ws.Rows.Item(33).Insert()
$ws.Rows.Item(33).RowHeight = 15

# Fill in the new case on row 33.
$ws.Cells.Item(33, 1).Value = "A 52224-2024"   # Beteckning
$ws.Cells.Item(33, 2).Value = 45607            # Datum (2024-11-11)
$ws.Cells.Item(33, 3).Value = 45623            # Förändrad (2024-11-27)
$ws.Cells.Item(33, 4).Value = "OKÄNT"          # Län
$ws.Cells.Item(33, 5).Value = "OKÄNT"          # Kommun
$ws.Cells.Item(33, 7).Value = 0.9              # Area (ha)
$ws.Cells.Item(33, 8).Value = 0                # Fridlysta
$ws.Cells.Item(33, 9).Value = 0                # Signalarter
$ws.Cells.Item(33, 10).Value = 0               # NT
$ws.Cells.Item(33, 11).Value = 0               # VU
$ws.Cells.Item(33, 12).Value = 0               # EN
$ws.Cells.Item(33, 13).Value = 0               # CR
$ws.Cells.Item(33, 14).Value = 0               # RE
$ws.Cells.Item(33, 15).Value = 0               # Rödlistade
$ws.Cells.Item(33, 16).Value = 0               # Hotade
$ws.Cells.Item(33, 17).Value = 0               # Alla arter

# The old row (now row 34, "A 52574-2024") also gets its "Förändrad"
# date refreshed to the new value.
$ws.Cells.Item(34, 3).Value = 45623
